$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Overview sheet: the zh-cn / de-de summary columns (E, F) for the two files
# that just got handed back move from "Ready for handoff" to the handed-back
# status. (Column G - last HO xliff generate date - is untouched.)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E4").Value = $newStatus
$wsOverview.Range("F4").Value = $newStatus
$wsOverview.Range("E5").Value = $newStatus
$wsOverview.Range("F5").Value = $newStatus

# ---------------------------------------------------------------------------
# zh-cn sheet: rows 4 (a489ea9e...) and 5 (d71a1c7d...) just got handed back.
#   - Status (col C) -> handed back
#   - Latest Target File (col I) now has the generated target file, linking
#     to the same source-file URL as column A
#   - Latest Handback File (col J) gets the xliff file name (same as col G)
#   - Latest Handback DateTime (col K) gets a real timestamp
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$zhLinks = @()
foreach ($h in $wsZhCn.Hyperlinks) { $zhLinks += $h }

$wsZhCn.Range("C4").Value = $newStatus
$wsZhCn.Range("I4").Value = "a489ea9e-d2b8-4b46-b298-146076d49ec4.md"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I4"), $zhLinks[4].Address, "", "", "a489ea9e-d2b8-4b46-b298-146076d49ec4.md")
$wsZhCn.Range("J4").Value = $wsZhCn.Range("G4").Value2
$wsZhCn.Range("K4").Value = "2016-10-14 09:10:34"

$wsZhCn.Range("C5").Value = $newStatus
$wsZhCn.Range("I5").Value = "d71a1c7d-5c20-4fbc-8e88-f8fc7c3b2692.md"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I5"), $zhLinks[5].Address, "", "", "d71a1c7d-5c20-4fbc-8e88-f8fc7c3b2692.md")
$wsZhCn.Range("J5").Value = $wsZhCn.Range("G5").Value2
$wsZhCn.Range("K5").Value = "2016-10-14 09:10:34"

# ---------------------------------------------------------------------------
# de-de sheet: same treatment.
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$deLinks = @()
foreach ($h in $wsDeDe.Hyperlinks) { $deLinks += $h }

$wsDeDe.Range("C4").Value = $newStatus
$wsDeDe.Range("I4").Value = "a489ea9e-d2b8-4b46-b298-146076d49ec4.md"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I4"), $deLinks[4].Address, "", "", "a489ea9e-d2b8-4b46-b298-146076d49ec4.md")
$wsDeDe.Range("J4").Value = $wsDeDe.Range("G4").Value2
$wsDeDe.Range("K4").Value = "2016-10-14 09:10:51"

$wsDeDe.Range("C5").Value = $newStatus
$wsDeDe.Range("I5").Value = "d71a1c7d-5c20-4fbc-8e88-f8fc7c3b2692.md"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I5"), $deLinks[5].Address, "", "", "d71a1c7d-5c20-4fbc-8e88-f8fc7c3b2692.md")
$wsDeDe.Range("J5").Value = $wsDeDe.Range("G5").Value2
$wsDeDe.Range("K5").Value = "2016-10-14 09:10:51"
